$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "passwprd" typo in the B1 header cell
$ws.Range("B1").Value = "password"

# Correct the email address used in A2 (also updates its hyperlink display text)
$ws.Range("A2").Value = "ntwano14@gmail.com"

# Give B2 a hyperlink (mirroring A2's) and fill in the password value,
# formatted with the built-in Hyperlink cell style
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:ntwano14@gmail.com")
$ws.Range("B2").Value = "NTWAno@16"
$ws.Range("B2").Style = "Hyperlink"

# New row 3 below, with A3 carrying the same Hyperlink style
$ws.Range("A3").Style = "Hyperlink"

# Leave the selection on B3, matching the final cursor position
$ws.Range("B3").Select() | Out-Null
